# Migrate validation to pandera
# Normalize the boolean-like strings in the "IsActive" column of the
# "data" sheet to a consistent "Yes"/"No" casing (the sheet previously
# had a stray lower-case "yes"/"no" alongside the canonical "Yes").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Row 4 ("Glass") had "yes" -> normalize to "Yes"
$ws.Range("G4").Value = "Yes"
# Row 5 ("Airbus A320") had "no" -> normalize to "No"
$ws.Range("G5").Value = "No"

# Move the active selection on the data sheet to match the saved view state.
[void]$ws.Activate()
[void]$ws.Range("I14").Select()
